$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "first"/"last" name columns are merged into a single "Full Name"
# column. Delete the old "last" column (C); cells in D:G shift left to
# become C:F, carrying their styles (e.g. the date style on the dob cell)
# with them.
$ws.Columns.Item(3).Delete()

# Restore the (now merged) column's width to match the old "dob" column
# that slid into slot C (23.2 chars wide).
$ws.Columns.Item(3).Width = 23.2

# Rename the header for the merged name column.
$ws.Range("B1").Value = "Full Name"

# Updated DB record.
$ws.Range("A2").Value = 111111111111
$ws.Range("B2").Value = "jhjh"
$ws.Range("C2").Value = 37336
$ws.Range("D2").Value = "F"
$ws.Range("E2").Value = "Rajasthan"

# Update the selection left over from editing the sheet.
$ws.Range("A3:E3").Select() | Out-Null
